$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at E. This shifts the existing Project/Fase/Code
# headers from E/F/G to F/G/H, carrying their column widths with them.
$ws.Columns("E").Insert()

# New header cells, in the order the author typed them: the two trailing
# columns (Medewerker, Tijd-ID) first, then the inserted column's header
# (MedewekerCode).
$ws.Range("I1").Value = "Medewerker"
$ws.Range("J1").Value = "Tijd-ID"
$ws.Range("E1").Value = "MedewekerCode"

# Column width tweaks for the two "new" text columns (author used
# AutoFit / bestFit in Excel; closest achievable custom widths below).
$ws.Columns("E").ColumnWidth = 15.666666666666666
$ws.Columns("I").ColumnWidth = 11.666666666666666

# Leave the selection on the newly inserted column, matching the saved
# workbook (an entire-column selection of E).
$ws.Range("E1:E1048576").Select() | Out-Null

# Defined (named) ranges added for the sync script to consume.
$wb.Names.Add('Codes', '=Sheet1!$H:$H') | Out-Null
$wb.Names.Add('Dates', '=Sheet1!$B:$B') | Out-Null
$wb.Names.Add('EmployeeCodes', '=Sheet1!$E:$E') | Out-Null
$wb.Names.Add('Employees', '=Sheet1!$I:$I') | Out-Null
$wb.Names.Add('IDs', '=Sheet1!$A:$A') | Out-Null
$wb.Names.Add('Phases', '=Sheet1!$G:$G') | Out-Null
$wb.Names.Add('Projects', '=Sheet1!$F:$F') | Out-Null
$wb.Names.Add('TimeBegin', '=Sheet1!$C:$C') | Out-Null
$wb.Names.Add('TimeEnd', '=Sheet1!$D:$D') | Out-Null
$wb.Names.Add('TimeIDs', '=Sheet1!$J:$J') | Out-Null
$wb.Names.Add('ID', '=Sheet1!#REF!') | Out-Null
